$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "E2" = 1
    "F2" = 0.3333333333333333
    "G2" = 0.04312
    "H2" = 0.12936
    "I2" = 0.06332032271043876
    "J2" = 0.06332032271043876
    "K2" = 1
    "L2" = 0.3333333333333333
    "M2" = 0.02747533333333333
    "N2" = 0.082426
    "O2" = 0.1853602590626968
    "P2" = 0.1853602590626968
    "Q2" = 0.001184736373333333
    "R2" = 0.01066262736
    "S2" = 0.01173707142154049
    "T2" = 0.01173707142154049
    "E3" = 1
    "F3" = 0.3333333333333333
    "G3" = 0.04312
    "H3" = 0.12936
    "I3" = 0.06332032271043876
    "J3" = 0.06332032271043876
    "O3" = 0.3818296302959431
    "P3" = 0.3818296302959431
    "Q3" = 0.002440477013333333
    "R3" = 0.02196429312
    "S3" = 0.02417757541074664
    "T3" = 0.02417757541074664
    "E4" = 1
    "F4" = 0.3333333333333333
    "G4" = 0.04312
    "H4" = 0.12936
    "I4" = 0.06332032271043876
    "J4" = 0.06332032271043876
    "O4" = 0.4328101106413601
    "P4" = 0.4328101106413601
    "Q4" = 0.00276632048
    "R4" = 0.02489688432
    "S4" = 0.02740567587815163
    "T4" = 0.02740567587815163
    "I5" = 0.5780441577995699
    "J5" = 0.5780441577995699
    "K5" = 1
    "L5" = 0.3333333333333333
    "M5" = 0.02747533333333333
    "N5" = 0.082426
    "O5" = 0.1853602590626968
    "P5" = 0.1853602590626968
    "Q5" = 0.01081532610422222
    "R5" = 0.09733793493799998
    "S5" = 0.1071464148394066
    "T5" = 0.1071464148394066
    "I6" = 0.5780441577995699
    "J6" = 0.5780441577995699
    "O6" = 0.3818296302959431
    "P6" = 0.3818296302959431
    "S6" = 0.2207143870673396
    "T6" = 0.2207143870673396
    "I7" = 0.5780441577995699
    "J7" = 0.5780441577995699
    "O7" = 0.4328101106413601
    "P7" = 0.4328101106413601
    "S7" = 0.2501833558928236
    "T7" = 0.2501833558928236
    "H8" = 0.7326729999999999
    "I8" = 0.3586355194899915
    "J8" = 0.3586355194899914
    "K8" = 1
    "L8" = 0.3333333333333333
    "M8" = 0.02747533333333333
    "N8" = 0.082426
    "O8" = 0.1853602590626968
    "P8" = 0.1853602590626968
    "Q8" = 0.006710144966444444
    "R8" = 0.06039130469799999
    "S8" = 0.06647677280174966
    "T8" = 0.06647677280174964
    "H9" = 0.7326729999999999
    "I9" = 0.3586355194899915
    "J9" = 0.3586355194899914
    "O9" = 0.3818296302959431
    "P9" = 0.3818296302959431
    "S9" = 0.1369376678178569
    "T9" = 0.1369376678178569
    "H10" = 0.7326729999999999
    "I10" = 0.3586355194899915
    "J10" = 0.3586355194899914
    "O10" = 0.4328101106413601
    "P10" = 0.4328101106413601
    "S10" = 0.1552210788703849
    "T10" = 0.1552210788703848
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
